$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Row 2 (005055865 / G3C / 115472.5) becomes 004328934 / VALERIA / 84000.
#    Force column A to text first so the leading zeros of the account number
#    survive the round-trip instead of being parsed as a plain number.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "004328934"
$ws.Cells.Item(2, 2).Value = "VALERIA"
$ws.Cells.Item(2, 3).Value = 84000

# 2. Insert a brand new row above row 4 (ANDRE) for 004486497 / ELENA / 12489.93
$ws.Rows.Item(4).Insert()
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "004486497"
$ws.Cells.Item(4, 2).Value = "ELENA"
$ws.Cells.Item(4, 3).Value = 12489.93

# 3. Remove the three now-trailing rows that used to sit right after ANDRE:
#    NATALIA/4000, ELENA/2503.77 and VALERIA/2000 (now at rows 6, 7 and 8
#    after the insertion above shifted everything down by one row).
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()
